# Apply the "Aufwand_git" worksheet update:
#  - add two new time-tracking rows (48 and 49) to "Tabelle1"
#  - these feed the existing SUM/ratio formulas in G2:G5, which recalc automatically
#  - update the sheet view's selection to the new last cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")  # "Tabelle1" is also $wb.ActiveSheet

# --- Row 48: 2024-03-05, 2h, "Webapp comparison for detail view property shapes"
$ws.Range("A48").NumberFormat = $ws.Range("A47").NumberFormat
$ws.Range("A48").Value = 45356
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = "Webapp comparison for detail view property shapes"

# --- Row 49: 2024-03-05, 3h, "Propsoal"
$ws.Range("A49").NumberFormat = $ws.Range("A47").NumberFormat
$ws.Range("A49").Value = 45356
$ws.Range("B49").Value = 3
$ws.Range("C49").Value = "Propsoal"

# Make sure the dependent summary formulas (SUM, ratio, etc.) are refreshed
$excel.CalculateFull()

# Move the view / selection to reflect the newly added last row, like in the
# saved workbook (scroll so the new rows are visible, select K49)
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K49").Select()
